$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $text) {
    $scratch = $ws.Range("Z1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.ClearContents()
}

Set-TextValue $ws "D2" "307.26"
Set-TextValue $ws "E2" "1.20%"
Set-TextValue $ws "E3" "1.36%"
Set-TextValue $ws "D4" "5.105"
Set-TextValue $ws "E4" "0.28%"
Set-TextValue $ws "D5" "0.08084"
Set-TextValue $ws "E5" "0.98%"
Set-TextValue $ws "D6" "1.942"
Set-TextValue $ws "E6" "0.60%"
Set-TextValue $ws "D7" "4.195"
Set-TextValue $ws "E7" "3.61%"
Set-TextValue $ws "D8" "7.747"
Set-TextValue $ws "E8" "-0.11%"
Set-TextValue $ws "D9" "0.9279"
Set-TextValue $ws "E9" "0.74%"
Set-TextValue $ws "D10" "0.1386"
Set-TextValue $ws "E10" "12.59%"
Set-TextValue $ws "D11" "0.1903"
Set-TextValue $ws "E11" "2.57%"
Set-TextValue $ws "D12" "0.09215"
Set-TextValue $ws "E12" "-2.68%"
Set-TextValue $ws "E13" "-5.02%"
Set-TextValue $ws "D14" "0.09831"
Set-TextValue $ws "E14" "-0.14%"
Set-TextValue $ws "D15" "0.001442"
Set-TextValue $ws "E15" "2.83%"
Set-TextValue $ws "D16" "0.005795"
Set-TextValue $ws "E16" "1.38%"
Set-TextValue $ws "D17" "3.617"
Set-TextValue $ws "E17" "3.64%"
Set-TextValue $ws "E18" "2.14%"
Set-TextValue $ws "D19" "0.3448"
Set-TextValue $ws "E19" "1.16%"
Set-TextValue $ws "E20" "4.44%"
Set-TextValue $ws "D21" "4.899"
Set-TextValue $ws "E21" "-2.99%"
Set-TextValue $ws "D23" "0.04432"
Set-TextValue $ws "E23" "-2.08%"
Set-TextValue $ws "D24" "0.001224"
Set-TextValue $ws "E24" "0.75%"
Set-TextValue $ws "D25" "0.004826"
Set-TextValue $ws "E25" "-0.45%"
Set-TextValue $ws "D26" "0.0001242"
Set-TextValue $ws "E26" "-0.63%"
Set-TextValue $ws "D39" "0.02029"
Set-TextValue $ws "E39" "5.05%"
Set-TextValue $ws "D40" "0.04919"
Set-TextValue $ws "E40" "3.41%"
Set-TextValue $ws "D41" "0.007610"
Set-TextValue $ws "E41" "1.25%"
Set-TextValue $ws "D42" "0.01008"
Set-TextValue $ws "E42" "5.46%"
Set-TextValue $ws "E43" "3.45%"
Set-TextValue $ws "D44" "0.002104"
Set-TextValue $ws "E44" "-0.30%"
Set-TextValue $ws "D45" "0.01101"
Set-TextValue $ws "E45" "0.09%"
Set-TextValue $ws "D46" "0.00006454"
Set-TextValue $ws "E46" "2.64%"
Set-TextValue $ws "E47" "0.18%"
Set-TextValue $ws "D48" "63.57"
Set-TextValue $ws "E48" "-1.41%"
Set-TextValue $ws "E49" "-19.85%"
Set-TextValue $ws "E50" "0.18%"
Set-TextValue $ws "E51" "0.18%"

$excel.CutCopyMode = $false
Write-Output "done"
